$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 85 (this pushes the old row 85 -> 86,
# old row 86 "Binary Search (basic)" header -> 87, etc. down by one),
# inheriting the formatting of the surrounding data rows.
$ws.Rows("85:85").Insert() | Out-Null

# Fill in the new entry: 111. Minimum Depth of Binary Tree
$ws.Range("A85").Value = 111
$ws.Range("B85").Value = "Minimum Depth of Binary Tree"
$ws.Range("C85").Value = "Easy"
$ws.Range("D85").Value = "DFS .recursion"
$ws.Range("E85").Value = 45858
$ws.Range("F85").Value = "Python"

# The "127. Word Ladder" hyperlink was anchored on the old B126; after the
# row insert that row is now B127, but hyperlink anchors don't auto-shift,
# so re-anchor it manually.
$ws.Range("B126").Hyperlinks.Delete() | Out-Null
$null = $ws.Hyperlinks.Add($ws.Range("B127"), "https://leetcode.com/problems/word-ladder/", "", "https://leetcode.com/problems/word-ladder/", "127. Word Ladder")

# Hyperlinks.Add re-stamps the cell with the built-in blue/underlined
# "Hyperlink" style; restore the sheet's own custom link look (10pt,
# no underline, #0A84FF) that the cell had before.
$f = $ws.Range("B127").Font
$f.Underline = -4142
$f.Size = 10
$f.Color = 16745482

# Match the selection shown in the author's edit (active cell moved to B85).
$ws.Range("B85").Select() | Out-Null
